# Update "想去人数" (F column) counts on the "展览" and "全部类型" sheets
# to reflect the newer scrape snapshot (output generated at 456a3b4).

$wb = $excel.ActiveWorkbook

# Map of worksheet name -> { row => new F value }
$updates = @{
    "展览" = @{
        2  = 166
        4  = 577
        5  = 1812
        9  = 2300
        10 = 118
        13 = 1406
        16 = 310
        24 = 76
        26 = 1435
        28 = 366
        29 = 157
        32 = 357
    }
    "全部类型" = @{
        2  = 166
        4  = 577
        5  = 1812
        10 = 2300
        11 = 118
        14 = 1406
        17 = 310
        25 = 76
        27 = 1435
        29 = 366
        30 = 157
        33 = 357
    }
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $rows = $updates[$sheetName]
    foreach ($row in $rows.Keys) {
        $newValue = $rows[$row]
        $ws.Range("F$row").Value = $newValue
    }
}
